# "Changed year for Monika's Chart and Taralynn"
# Each per-year sheet gets renamed to match the new year label, the year
# value shown in A2 is bumped by one, and the sheet's active-cell selection
# is moved to where the editor last clicked.

$wb = $excel.ActiveWorkbook

$sheetInfo = @(
    @{ Index = 1; Name = "2016"; Selection = "E10" },
    @{ Index = 2; Name = "2017"; Selection = "O20" },
    @{ Index = 3; Name = "2018"; Selection = "A3"  },
    @{ Index = 4; Name = "2019"; Selection = "A3"  },
    @{ Index = 5; Name = "2020"; Selection = "I14" }
)

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Index)

    # Rename the sheet tab (Sheet1..Sheet5 -> 2016..2020)
    $ws.Name = $info.Name

    # Bump the displayed model year by one (2015->2016, ..., 2019->2020)
    $yearCell = $ws.Range("A2")
    $currentYear = $yearCell.Value()
    $yearCell.Value = $currentYear + 1

    # Move the active cell / selection to where editing last left off
    [void]$ws.Range($info.Selection).Select()
}

# Keep the last sheet ("2020") as the active/visible tab, same as before
[void]$wb.Worksheets.Item(5).Activate()
